$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(800, 100, 0.5, 0.6, 7123, 0),
    @(800, 50,  0.5, 0.6, 3395, 0),
    @(800, 50,  0.5, 0.6, 3409, 0),
    @(800, 50,  0.5, 0.6, 3344, 0),
    @(800, 50,  0.5, 0.6, 3089, 0),
    @(800, 50,  0.5, 0.6, 2107, 0)
)

$startRow = 229
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
